$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "done" list in column B down by one starting at row 37,
# then add the two new entries, matching the commit
# "manual input an ongoing timing record".

$ws.Cells.Item(37, 2).Value = "Refresh last activity timer after last record is deleted"
$ws.Cells.Item(37, 2).Style = "Good"

$ws.Cells.Item(38, 2).Value = "title total amout 每日自動清零 (nappy)"
$ws.Cells.Item(38, 2).Style = "Good"

$ws.Cells.Item(39, 2).Value = "title total amout 每日自動清零後 reset title"
$ws.Cells.Item(39, 2).Style = "Good"

$ws.Cells.Item(40, 2).Value = "manual input an ongoing timing record"
$ws.Cells.Item(40, 2).Style = "Good"

$ws.Cells.Item(41, 2).Value = "finish an ongoing timing record before delete it"
$ws.Cells.Item(41, 2).Style = "Good"

$ws.Cells.Item(42, 2).ClearContents() | Out-Null

# Update the view: active window selection
$ws.Activate() | Out-Null
$ws.Range("B43").Select() | Out-Null
